$d = $word.ActiveDocument

# --- Locate the "KEY ACHIEVEMENTS AND IMPACT" section bounds -----------------
# (the section runs from that Heading2 paragraph up to, but excluding, the
#  next Heading2 paragraph - "TECHNICAL SKILLS"). Several bullet strings in
# this section are exact duplicates of text used elsewhere in the resume
# (e.g. the "Trigonometric algorithm..." bullet also appears verbatim under
# "PROFESSIONAL EXPERIENCE"), so all Find/Replace calls below are scoped to
# this section's Range only - never to the whole document - to avoid
# touching the other occurrences.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($startPara -eq $null -and $p.Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $startPara = $i
        continue
    }
    if ($startPara -ne $null -and $p.Range.ParagraphStyle.NameLocal -eq "Heading 2") {
        $endPara = $i
        break
    }
}

function Get-SectionRange($doc, $first, $last) {
    $s = $doc.Paragraphs.Item($first).Range.Start
    $e = $doc.Paragraphs.Item($last).Range.Start
    return $doc.Range($s, $e)
}

# --- Rewrite the bullets as impact-focused accomplishment statements ---------

# 1) "Discovered systematic race coding errors..." -> algorithmic-innovation bullet
$rng = Get-SectionRange $d $startPara $endPara
$rng.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%",
    2) | Out-Null

# 2) "Built redistricting platform..." -> "$4.7M savings enabled nonprofit access"
$rng = Get-SectionRange $d $startPara $endPara
$rng.Find.Execute(
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access",
    2) | Out-Null

# 3) "Trigonometric algorithm for boundary estimation..." -> demographic-discovery bullet
$rng = Get-SectionRange $d $startPara $endPara
$rng.Find.Execute(
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    2) | Out-Null

# 4) "Developed longitudinal data analysis methods..." -> accuracy-improvement bullet
$rng = Get-SectionRange $d $startPara $endPara
$rng.Find.Execute(
    "Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "178% accuracy improvement in racial classification algorithms",
    2) | Out-Null

# --- Drop the two remaining bullets (ETL pipelines / AWS data warehouse) -----
# Walk the section backwards so deleting a paragraph doesn't shift the index
# of paragraphs not yet visited.
for ($i = $endPara - 1; $i -ge $startPara; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets*" -or
        $t -like "*Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy*") {
        $p.Range.Delete()
    }
}
